# Apply the "ajustes de los datos suavizados" edits to the
# ECOContSmooth1min sheet: a handful of rows that held near-garbage
# smoothed turbidity statistics (columns B/C/D) are cleared out, and a
# few other rows get corrected Mean/Std/CV values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECOContSmooth1min")

# --- Rows whose turbidity Mean/Std/CV (B/C/D) need to be blanked out ---
$rowsToClear = 2, 3, 4, 168, 441, 445
foreach ($r in $rowsToClear) {
    $addr = "B" + $r + ":D" + $r
    $ws.Range($addr).ClearContents()
}

# Rows 442-444 only need column B and C cleared; column D is already blank.
$rowsToClearBC = 442, 443, 444
foreach ($r in $rowsToClearBC) {
    $addr = "B" + $r + ":C" + $r
    $ws.Range($addr).ClearContents()
}

# --- Rows whose turbidity Mean/Std/CV (B/C/D) get corrected values ---
$ws.Range("B5").Value = 25.40396
$ws.Range("C5").Value = 1.054276447617039
$ws.Range("D5").Value = 4.150047660353104

$ws.Range("B114").Value = 26.49293333333334
$ws.Range("C114").Value = 0.1407579956284309
$ws.Range("D114").Value = 0.5313039287021101

$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 0

$ws.Range("B239").Value = 25.9647
$ws.Range("C239").Value = 0.1723926332532804
$ws.Range("D239").Value = 0.6639500292831438

$ws.Range("B440").Value = 39.08926666666667
$ws.Range("C440").Value = 0.372410651476745
$ws.Range("D440").Value = 0.9527184396997598
